$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: MuSCs -> Calca -> Calcr -> ECs
$ws.Range("A2").Value = "MuSCs"
$ws.Range("B2").Value = "Calca"
$ws.Range("C2").Value = "Calcr"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.5
$ws.Range("G2").Value = 0.354751
$ws.Range("H2").Value = 0.709502
$ws.Range("I2").Value = 0.7067029850439027
$ws.Range("J2").Value = 0.6163205031315422
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.5
$ws.Range("M2").Value = 0.046576
$ws.Range("N2").Value = 0.093152
$ws.Range("O2").Value = 0.07994370161857847
$ws.Range("P2").Value = 0.07994370161857847
$ws.Range("Q2").Value = 0.016522882576
$ws.Range("R2").Value = 0.066091530304
$ws.Range("S2").Value = 0.05649645256930848
$ws.Range("T2").Value = 0.04927094240376017

# Row 3: MuSCs -> Calca -> Calcr -> MuSCs
$ws.Range("A3").Value = "MuSCs"
$ws.Range("B3").Value = "Calca"
$ws.Range("C3").Value = "Calcr"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.5
$ws.Range("G3").Value = 0.354751
$ws.Range("H3").Value = 0.709502
$ws.Range("I3").Value = 0.7067029850439027
$ws.Range("J3").Value = 0.6163205031315422
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.536034
$ws.Range("N3").Value = 1.072068
$ws.Range("O3").Value = 0.9200562983814217
$ws.Range("P3").Value = 0.9200562983814217
$ws.Range("Q3").Value = 0.190158597534
$ws.Range("R3").Value = 0.760634390136
$ws.Range("S3").Value = 0.6502065324745944
$ws.Range("T3").Value = 0.5670495607277821

# Row 4: Neutrophils -> Calca -> Calcr -> ECs
$ws.Range("A4").Value = "Neutrophils"
$ws.Range("B4").Value = "Calca"
$ws.Range("C4").Value = "Calcr"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1472293333333334
$ws.Range("H4").Value = 0.441688
$ws.Range("I4").Value = 0.2932970149560972
$ws.Range("J4").Value = 0.3836794968684579
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.5
$ws.Range("M4").Value = 0.046576
$ws.Range("N4").Value = 0.093152
$ws.Range("O4").Value = 0.07994370161857847
$ws.Range("P4").Value = 0.07994370161857847
$ws.Range("Q4").Value = 0.006857353429333334
$ws.Range("R4").Value = 0.041144120576
$ws.Range("S4").Value = 0.02344724904926998
$ws.Range("T4").Value = 0.03067275921481831

# Row 5: Neutrophils -> Calca -> Calcr -> MuSCs
$ws.Range("A5").Value = "Neutrophils"
$ws.Range("B5").Value = "Calca"
$ws.Range("C5").Value = "Calcr"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.1472293333333334
$ws.Range("H5").Value = 0.441688
$ws.Range("I5").Value = 0.2932970149560972
$ws.Range("J5").Value = 0.3836794968684579
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.536034
$ws.Range("N5").Value = 1.072068
$ws.Range("O5").Value = 0.9200562983814217
$ws.Range("P5").Value = 0.9200562983814217
$ws.Range("Q5").Value = 0.07891992846400001
$ws.Range("R5").Value = 0.473519570784
$ws.Range("S5").Value = 0.2698497659068273
$ws.Range("T5").Value = 0.3530067376536397
